# "make email and evening phone optional"
#
# The sheet stores a flat key/value config table in columns A/B. The key
# "form_another" (row 16) holds the HTML message shown to the parent after a
# successful submission. Since email is now optional, the copy can no longer
# promise "you will recieve an email shortly" - drop that sentence while
# keeping the rest of the message (and its <br> formatting) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetKey = "form_another"
$oldValue  = "Your submission was successful!<br>You will recieve an email shortly.<br><br>Do you wish to submit again for another child?"
$newValue  = "Your submission was successful!<br><br>Do you wish to submit again for another child?"

# Locate the row holding the key in column A so the edit is resilient to the
# sheet being reordered, rather than hard-coding row 16.
$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($i = 1; $i -le $rowCount; $i++) {
    $key = $ws.Cells.Item($i, 1).Text
    if ($key -eq $targetKey) {
        $cell = $ws.Cells.Item($i, 2)
        if ($cell.Text -eq $oldValue) {
            $cell.Value = $newValue
        }
        break
    }
}
